$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column E metadata rows describe the "provincia" field.
# It is being re-curated from an sdmx dimension (refArea) to an iaest measure.
$ws.Range("E2").Value = "iaest-measure:provincia"
$ws.Range("E3").Value = "medida"
$ws.Range("E4").Value = "xsd:int"
